$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(62, 8).Value = 4997  # H62
$ws.Cells.Item(62, 10).Value = 4994  # J62
$ws.Cells.Item(62, 12).Value = 4994  # L62
$ws.Cells.Item(62, 14).Value = -6242  # N62
$ws.Cells.Item(65, 8).Value = 4997  # H65
$ws.Cells.Item(65, 10).Value = 4994  # J65
$ws.Cells.Item(65, 12).Value = 24970  # L65
$ws.Cells.Item(65, 14).Value = -31210  # N65
$ws.Cells.Item(103, 8).Value = 4000.25  # H103
$ws.Cells.Item(111, 8).Value = 3059.7273  # H111
$ws.Cells.Item(111, 9).Value = 1047.5  # I111
$ws.Cells.Item(111, 10).Value = 3506.889  # J111
$ws.Cells.Item(111, 11).Value = 3142.5  # K111
$ws.Cells.Item(111, 12).Value = 10520.667  # L111
$ws.Cells.Item(111, 13).Value = -75.5  # M111
$ws.Cells.Item(111, 14).Value = -16654.667  # N111

# Sheet: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(42, 8).Value = 0  # H42
$ws.Cells.Item(42, 9).Value = 0  # I42
$ws.Cells.Item(42, 11).Value = 0  # K42
$ws.Cells.Item(42, 13).ClearContents()  # M42
$ws.Cells.Item(61, 8).Value = 11583.333  # H61
$ws.Cells.Item(61, 9).Value = 7875  # I61
$ws.Cells.Item(61, 11).Value = 7875  # K61
$ws.Cells.Item(61, 13).Value = -7663  # M61
$ws.Cells.Item(136, 8).Value = 11583.333  # H136
$ws.Cells.Item(136, 9).Value = 7875  # I136
$ws.Cells.Item(136, 11).Value = 23625  # K136
$ws.Cells.Item(136, 13).Value = -21075  # M136

# Sheet: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(20, 8).Value = 1184.8334  # H20
$ws.Cells.Item(20, 9).Value = 1249.75  # I20
$ws.Cells.Item(20, 11).Value = 1249.75  # K20
$ws.Cells.Item(20, 13).Value = -1002.75  # M20

# Sheet: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(31, 8).Value = 12400  # H31
$ws.Cells.Item(31, 9).Value = 9000  # I31
$ws.Cells.Item(31, 11).Value = 9000  # K31
$ws.Cells.Item(31, 13).Value = -8705  # M31
$ws.Cells.Item(34, 8).Value = 12400  # H34
$ws.Cells.Item(34, 9).Value = 9000  # I34
$ws.Cells.Item(34, 11).Value = 9000  # K34
$ws.Cells.Item(34, 13).Value = -8798  # M34
$ws.Cells.Item(58, 8).Value = 5486  # H58
$ws.Cells.Item(58, 9).Value = 3314.7778  # I58
$ws.Cells.Item(58, 11).Value = 3314.7778  # K58
$ws.Cells.Item(58, 13).Value = -3111.7778  # M58
$ws.Cells.Item(98, 8).Value = 0  # H98
$ws.Cells.Item(98, 10).Value = 0  # J98
$ws.Cells.Item(98, 12).Value = 0  # L98
$ws.Cells.Item(98, 14).ClearContents()  # N98
$ws.Cells.Item(122, 8).Value = 1349.5  # H122
$ws.Cells.Item(122, 9).Value = 1200  # I122
$ws.Cells.Item(122, 10).Value = 1499  # J122
$ws.Cells.Item(122, 11).Value = 3600  # K122
$ws.Cells.Item(122, 12).Value = 4497  # L122
$ws.Cells.Item(122, 13).Value = -1150  # M122
$ws.Cells.Item(122, 14).Value = -9397  # N122
$ws.Cells.Item(136, 8).Value = 5486  # H136
$ws.Cells.Item(136, 9).Value = 3314.7778  # I136
$ws.Cells.Item(136, 11).Value = 9944.3334  # K136
$ws.Cells.Item(136, 13).Value = -7394.3334  # M136

# Sheet: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(33, 8).Value = 250043  # H33
$ws.Cells.Item(33, 10).Value = 333371  # J33
$ws.Cells.Item(33, 12).Value = 2000226  # L33
$ws.Cells.Item(33, 14).Value = -2000792  # N33
$ws.Cells.Item(38, 8).Value = 65.59999999999999  # H38
$ws.Cells.Item(38, 9).Value = 85  # I38
$ws.Cells.Item(38, 10).Value = 52.666668  # J38
$ws.Cells.Item(38, 11).Value = 255  # K38
$ws.Cells.Item(38, 12).Value = 158.000004  # L38
$ws.Cells.Item(38, 13).Value = 92  # M38
$ws.Cells.Item(38, 14).Value = -852.000004  # N38
$ws.Cells.Item(50, 8).Value = 12.5  # H50
$ws.Cells.Item(50, 9).Value = 5  # I50
$ws.Cells.Item(50, 10).Value = 20  # J50
$ws.Cells.Item(50, 11).Value = 15  # K50
$ws.Cells.Item(50, 12).Value = 60  # L50
$ws.Cells.Item(50, 13).Value = 466  # M50
$ws.Cells.Item(50, 14).Value = -1022  # N50
$ws.Cells.Item(51, 8).Value = 2200  # H51
$ws.Cells.Item(51, 10).Value = 2200  # J51
$ws.Cells.Item(51, 12).Value = 6600  # L51
$ws.Cells.Item(51, 14).Value = -7520  # N51
$ws.Cells.Item(53, 8).Value = 12.5  # H53
$ws.Cells.Item(53, 9).Value = 5  # I53
$ws.Cells.Item(53, 10).Value = 20  # J53
$ws.Cells.Item(53, 11).Value = 15  # K53
$ws.Cells.Item(53, 12).Value = 60  # L53
$ws.Cells.Item(53, 13).Value = 466  # M53
$ws.Cells.Item(53, 14).Value = -1022  # N53
$ws.Cells.Item(104, 8).Value = 2999  # H104
$ws.Cells.Item(104, 10).Value = 2999  # J104
$ws.Cells.Item(104, 12).Value = 8997  # L104
$ws.Cells.Item(104, 14).Value = -14239  # N104
$ws.Cells.Item(119, 8).Value = 0  # H119
$ws.Cells.Item(119, 9).Value = 0  # I119
$ws.Cells.Item(119, 11).Value = 0  # K119
$ws.Cells.Item(119, 13).ClearContents()  # M119
$ws.Cells.Item(131, 8).Value = 2161.2222  # H131
$ws.Cells.Item(131, 9).Value = 943.5714  # I131
$ws.Cells.Item(131, 11).Value = 2830.7142  # K131
$ws.Cells.Item(131, 13).Value = 2209.2858  # M131

# Sheet: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(2, 8).Value = 51.57143  # H2
$ws.Cells.Item(2, 9).Value = 47.846153  # I2
$ws.Cells.Item(2, 10).Value = 100  # J2
$ws.Cells.Item(2, 11).Value = 47.846153  # K2
$ws.Cells.Item(2, 12).Value = 100  # L2
$ws.Cells.Item(2, 13).Value = 65.153847  # M2
$ws.Cells.Item(2, 14).Value = -326  # N2
$ws.Cells.Item(126, 8).Value = 5900  # H126
$ws.Cells.Item(126, 10).Value = 6666.6665  # J126
$ws.Cells.Item(126, 12).Value = 19999.9995  # L126
$ws.Cells.Item(126, 14).Value = -24939.9995  # N126

# Sheet: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(21, 8).Value = 0  # H21
$ws.Cells.Item(21, 10).Value = 0  # J21
$ws.Cells.Item(21, 12).Value = 0  # L21
$ws.Cells.Item(21, 14).ClearContents()  # N21
$ws.Cells.Item(22, 8).Value = 0  # H22
$ws.Cells.Item(22, 9).Value = 0  # I22
$ws.Cells.Item(22, 11).Value = 0  # K22
$ws.Cells.Item(22, 13).ClearContents()  # M22
$ws.Cells.Item(27, 8).Value = 0  # H27
$ws.Cells.Item(27, 9).Value = 0  # I27
$ws.Cells.Item(27, 11).Value = 0  # K27
$ws.Cells.Item(27, 13).ClearContents()  # M27
$ws.Cells.Item(35, 8).Value = 2661.25  # H35
$ws.Cells.Item(35, 9).Value = 1898.5714  # I35
$ws.Cells.Item(35, 10).Value = 8000  # J35
$ws.Cells.Item(35, 11).Value = 1898.5714  # K35
$ws.Cells.Item(35, 12).Value = 8000  # L35
$ws.Cells.Item(35, 13).Value = -1562.5714  # M35
$ws.Cells.Item(35, 14).Value = -8672  # N35
$ws.Cells.Item(46, 8).Value = 5000.1665  # H46
$ws.Cells.Item(46, 9).Value = 3750.25  # I46
$ws.Cells.Item(46, 11).Value = 3750.25  # K46
$ws.Cells.Item(46, 13).Value = -3562.25  # M46
$ws.Cells.Item(55, 8).Value = 1537  # H55
$ws.Cells.Item(55, 9).Value = 2242.5  # I55
$ws.Cells.Item(55, 10).Value = 1066.6666  # J55
$ws.Cells.Item(55, 11).Value = 2242.5  # K55
$ws.Cells.Item(55, 12).Value = 1066.6666  # L55
$ws.Cells.Item(55, 13).Value = -2069.5  # M55
$ws.Cells.Item(55, 14).Value = -1412.6666  # N55
$ws.Cells.Item(82, 8).Value = 1666.6666  # H82
$ws.Cells.Item(82, 9).Value = 1666.6666  # I82
$ws.Cells.Item(82, 11).Value = 1666.6666  # K82
$ws.Cells.Item(82, 13).Value = -1305.6666  # M82
$ws.Cells.Item(85, 8).Value = 1666.6666  # H85
$ws.Cells.Item(85, 9).Value = 1666.6666  # I85
$ws.Cells.Item(85, 11).Value = 1666.6666  # K85
$ws.Cells.Item(85, 13).Value = -418.6666  # M85
$ws.Cells.Item(122, 8).Value = 4000.5  # H122
$ws.Cells.Item(122, 9).Value = 3004  # I122
$ws.Cells.Item(122, 11).Value = 9012  # K122
$ws.Cells.Item(122, 13).Value = -6562  # M122

# Sheet: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(54, 8).Value = 31500  # H54
$ws.Cells.Item(54, 10).Value = 31500  # J54
$ws.Cells.Item(54, 12).Value = 31500  # L54
$ws.Cells.Item(54, 14).Value = -32540  # N54
$ws.Cells.Item(55, 8).Value = 10000  # H55
$ws.Cells.Item(55, 9).Value = 0  # I55
$ws.Cells.Item(55, 10).Value = 10000  # J55
$ws.Cells.Item(55, 11).Value = 0  # K55
$ws.Cells.Item(55, 12).Value = 10000  # L55
$ws.Cells.Item(55, 13).ClearContents()  # M55
$ws.Cells.Item(55, 14).Value = -10554  # N55
$ws.Cells.Item(107, 8).Value = 635.5  # H107
$ws.Cells.Item(107, 9).Value = 471  # I107
$ws.Cells.Item(107, 10).Value = 800  # J107
$ws.Cells.Item(107, 11).Value = 1413  # K107
$ws.Cells.Item(107, 12).Value = 2400  # L107
$ws.Cells.Item(107, 13).Value = 507  # M107
$ws.Cells.Item(107, 14).Value = -6240  # N107
$ws.Cells.Item(125, 8).Value = 37150  # H125
$ws.Cells.Item(125, 10).Value = 37150  # J125
$ws.Cells.Item(125, 12).Value = 37150  # L125
$ws.Cells.Item(125, 14).Value = -46990  # N125
$ws.Cells.Item(136, 8).Value = 7999.4707  # H136
$ws.Cells.Item(136, 9).Value = 6213.643  # I136
$ws.Cells.Item(136, 11).Value = 18640.929  # K136
$ws.Cells.Item(136, 13).Value = -16090.929  # M136
